# Add technology rows for heat storages (discharger / charger) and
# flip the sheet's AutoFilter from "Tampere_dheat" to "Helsinki_dheat".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 222: Heat storage discharger (Heat_capa = 100)
$ws.Cells.Item(222, 1).Value = "FI00"
$ws.Cells.Item(222, 2).Value = "Helsinki_dheat"
$ws.Cells.Item(222, 3).Value = "Heat storage discharger"
$ws.Cells.Item(222, 4).Value = "Distributed Energy"
$ws.Cells.Item(222, 5).Value = 2040
$ws.Cells.Item(222, 7).Value = 100

# New row 223: Heat storage charger (Other_capa = 100)
$ws.Cells.Item(223, 1).Value = "FI00"
$ws.Cells.Item(223, 2).Value = "Helsinki_dheat"
$ws.Cells.Item(223, 3).Value = "Heat storage charger"
$ws.Cells.Item(223, 4).Value = "Distributed Energy"
$ws.Cells.Item(223, 5).Value = 2040
$ws.Cells.Item(223, 8).Value = 100

# Re-point the AutoFilter on column B (Heatnode) at "Helsinki_dheat"
# instead of "Tampere_dheat" -- this also recomputes which rows are
# hidden/shown for the whole filtered range.
$ws.Range("A1:J223").AutoFilter(2, @("Helsinki_dheat"), 7) | Out-Null

# Move the active selection to reflect where editing left off.
$ws.Range("H224").Select() | Out-Null
